$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 152480
$ws.Range("C4").Value = 144041
$ws.Range("C5").Value = 8439
$ws.Range("C7").Value = 5.53
$ws.Range("C8").Value = 63.76
